# fix: tolerate empty columns in tabular files
#
# The data previously lived in column B (with an extra formatting-only
# column C to its right). To tolerate "empty" leading columns in the
# source tabular files, two blank columns are inserted before column B,
# pushing the former column-B data into column D (and the extra
# formatting column from C into E/F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at B:C - this shifts the existing B (and
# beyond) content two columns to the right, carrying formatting along.
$ws.Range("B1:C1").EntireColumn.Insert()

# Reflect the author's new cursor position after the edit.
$ws.Range("I17").Select()
